$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 9 new rows right after the header row (row 1). This pushes the
# existing data rows (old rows 2-21) down to rows 11-30, matching the diff.
$ws.Range("A2:A10").EntireRow.Insert()

# Row-insert copies formatting from the row above (the bold header row),
# which would stamp a style index on the new cells. The target diff shows
# these new data cells with no style attribute (like the rest of the data
# rows), so clear any inherited formatting.
$ws.Range("A2:C10").ClearFormats()

# Values for the newly inserted rows 2-10.
$newData = @(
  @(-0.007375299738829639, -0.04188020327402392, 0.02404832670136416),
  @(0.009710959871025781, -0.04103577224647274, -0.004662338863401103),
  @(0.04170951860792497, 0.04793495536946196, -0.07793023174299907),
  @(-0.02235946409842538, -0.01261257046066655, 0.04359601744834121),
  @(-0.03681361302733401, -0.006297301829737623, -0.01943090470398167),
  @(0.02995036389021311, 0.01477754981640481, -0.02054483487325551),
  @(0.04258089907029088, -0.09914881779867053, -0.03582545123336939),
  @(0.0959058403968811, -0.1950187236070633, 0.0042760567739605),
  @(-0.04207783586838668, -0.3005187625394148, 0.05587620359352399)
)

$r = 2
foreach ($row in $newData) {
  $ws.Cells.Item($r, 1).Value = $row[0]
  $ws.Cells.Item($r, 2).Value = $row[1]
  $ws.Cells.Item($r, 3).Value = $row[2]
  $r++
}

# Append one more new row (row 31) after the (now shifted) last data row (30).
$ws.Cells.Item(31, 1).Value = -0.009441461181268051
$ws.Cells.Item(31, 2).Value = 0.03861925794797785
$ws.Cells.Item(31, 3).Value = -0.04366788180435401
